$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "status" column entirely (column A), shifting every other
# column one place to the left.
$ws.Range("A1").EntireColumn.Delete()

# Clear the related_resource value for the first data row (now row 2,
# column F) - it is no longer populated for this entry.
$ws.Range("F2").ClearContents()

# Move the active selection, matching the saved view state.
$ws.Range("E13").Select()
